$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old table first (rows 1-5, cols A-D from the previous layout)
$ws.Range("A1:D5").ClearContents()

# Row 2: data row (entered first, D2 left for later)
$ws.Range("A2").Value = "jeden"
$ws.Range("B2").Value = "dwa "
$ws.Range("C2").Value = "trzy"
$ws.Range("E2").Value = "XD2"
$ws.Range("F2").Value = "XD3"

# Row 3: data row
$ws.Range("A3").Value = "jeden"
$ws.Range("B3").Value = "dwa "
$ws.Range("C3").Value = "trzy"
$ws.Range("D3").Value = "cztery"
$ws.Range("E3").Value = "piec"
$ws.Range("F3").Value = "szesc"
$ws.Range("G3").Value = "siedem"

# Row 1: header/title row
$ws.Range("A1").Value = "Title1"
$ws.Range("B1").Value = "Title2"
$ws.Range("C1").Value = "title3"
$ws.Range("D1").Value = "title4"
$ws.Range("E1").Value = "title5"
$ws.Range("F1").Value = "title6"
$ws.Range("G1").Value = "title7"

# Finally fill in D2 with the long text, added last
$ws.Range("D2").Value = "Długi tekst, długi tekst, długi tekst Długi tekst, długi tekst, długi tekst"

$ws.Range("D2").Select()
